$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Hp"
$ws.Cells.Item(2,3).Value2 = "Itgam"
$ws.Cells.Item(2,4).Value2 = "MuSCs"
$ws.Cells.Item(2,5).Value2 = 1
$ws.Cells.Item(2,6).Value2 = 0.3333333333333333
$ws.Cells.Item(2,7).Value2 = 0.1238986666666667
$ws.Cells.Item(2,8).Value2 = 0.371696
$ws.Cells.Item(2,9).Value2 = 0.01923905185495286
$ws.Cells.Item(2,10).Value2 = 0.01923905185495286
$ws.Cells.Item(2,11).Value2 = 1
$ws.Cells.Item(2,12).Value2 = 0.3333333333333333
$ws.Cells.Item(2,13).Value2 = 0.02725333333333333
$ws.Cells.Item(2,14).Value2 = 0.08176
$ws.Cells.Item(2,15).Value2 = 0.0007089206372884383
$ws.Cells.Item(2,16).Value2 = 0.0007089206372884382
$ws.Cells.Item(2,17).Value2 = 0.003376651662222222
$ws.Cells.Item(2,18).Value2 = 0.03038986496
$ws.Cells.Item(2,19).Value2 = 0.00001363896090183849
$ws.Cells.Item(2,20).Value2 = 0.00001363896090183849

# Row 3
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Hp"
$ws.Cells.Item(3,3).Value2 = "Itgam"
$ws.Cells.Item(3,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(3,5).Value2 = 1
$ws.Cells.Item(3,6).Value2 = 0.3333333333333333
$ws.Cells.Item(3,7).Value2 = 0.1238986666666667
$ws.Cells.Item(3,8).Value2 = 0.371696
$ws.Cells.Item(3,9).Value2 = 0.01923905185495286
$ws.Cells.Item(3,10).Value2 = 0.01923905185495286
$ws.Cells.Item(3,11).Value2 = 3
$ws.Cells.Item(3,12).Value2 = 1
$ws.Cells.Item(3,13).Value2 = 38.416166
$ws.Cells.Item(3,14).Value2 = 115.248498
$ws.Cells.Item(3,15).Value2 = 0.9992910793627116
$ws.Cells.Item(3,16).Value2 = 0.9992910793627116
$ws.Cells.Item(3,17).Value2 = 4.759711745845333
$ws.Cells.Item(3,18).Value2 = 42.837405712608
$ws.Cells.Item(3,19).Value2 = 0.01922541289405102
$ws.Cells.Item(3,20).Value2 = 0.01922541289405102

# Row 4
$ws.Cells.Item(4,1).Value2 = "FAPs"
$ws.Cells.Item(4,2).Value2 = "Hp"
$ws.Cells.Item(4,3).Value2 = "Itgam"
$ws.Cells.Item(4,4).Value2 = "MuSCs"
$ws.Cells.Item(4,5).Value2 = 3
$ws.Cells.Item(4,6).Value2 = 1
$ws.Cells.Item(4,7).Value2 = 5.810518333333333
$ws.Cells.Item(4,8).Value2 = 17.431555
$ws.Cells.Item(4,9).Value2 = 0.9022604239955847
$ws.Cells.Item(4,10).Value2 = 0.9022604239955845
$ws.Cells.Item(4,11).Value2 = 1
$ws.Cells.Item(4,12).Value2 = 0.3333333333333333
$ws.Cells.Item(4,13).Value2 = 0.02725333333333333
$ws.Cells.Item(4,14).Value2 = 0.08176
$ws.Cells.Item(4,15).Value2 = 0.0007089206372884383
$ws.Cells.Item(4,16).Value2 = 0.0007089206372884382
$ws.Cells.Item(4,17).Value2 = 0.1583559929777778
$ws.Cells.Item(4,18).Value2 = 1.4252039368
$ws.Cells.Item(4,19).Value2 = 0.0006396310347790865
$ws.Cells.Item(4,20).Value2 = 0.0006396310347790862

# Row 5
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Hp"
$ws.Cells.Item(5,3).Value2 = "Itgam"
$ws.Cells.Item(5,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(5,5).Value2 = 3
$ws.Cells.Item(5,6).Value2 = 1
$ws.Cells.Item(5,7).Value2 = 5.810518333333333
$ws.Cells.Item(5,8).Value2 = 17.431555
$ws.Cells.Item(5,9).Value2 = 0.9022604239955847
$ws.Cells.Item(5,10).Value2 = 0.9022604239955845
$ws.Cells.Item(5,11).Value2 = 3
$ws.Cells.Item(5,12).Value2 = 1
$ws.Cells.Item(5,13).Value2 = 38.416166
$ws.Cells.Item(5,14).Value2 = 115.248498
$ws.Cells.Item(5,15).Value2 = 0.9992910793627116
$ws.Cells.Item(5,16).Value2 = 0.9992910793627116
$ws.Cells.Item(5,17).Value2 = 223.2178368393766
$ws.Cells.Item(5,18).Value2 = 2008.96053155439
$ws.Cells.Item(5,19).Value2 = 0.9016207929608057
$ws.Cells.Item(5,20).Value2 = 0.9016207929608054

# Row 6
$ws.Cells.Item(6,1).Value2 = "MuSCs"
$ws.Cells.Item(6,2).Value2 = "Hp"
$ws.Cells.Item(6,3).Value2 = "Itgam"
$ws.Cells.Item(6,4).Value2 = "MuSCs"
$ws.Cells.Item(6,5).Value2 = 1
$ws.Cells.Item(6,6).Value2 = 0.3333333333333333
$ws.Cells.Item(6,7).Value2 = 0.1062546666666667
$ws.Cells.Item(6,8).Value2 = 0.318764
$ws.Cells.Item(6,9).Value2 = 0.01649928200866351
$ws.Cells.Item(6,10).Value2 = 0.01649928200866351
$ws.Cells.Item(6,11).Value2 = 1
$ws.Cells.Item(6,12).Value2 = 0.3333333333333333
$ws.Cells.Item(6,13).Value2 = 0.02725333333333333
$ws.Cells.Item(6,14).Value2 = 0.08176
$ws.Cells.Item(6,15).Value2 = 0.0007089206372884383
$ws.Cells.Item(6,16).Value2 = 0.0007089206372884382
$ws.Cells.Item(6,17).Value2 = 0.002895793848888889
$ws.Cells.Item(6,18).Value2 = 0.02606214464
$ws.Cells.Item(6,19).Value2 = 0.0000116966815163834
$ws.Cells.Item(6,20).Value2 = 0.0000116966815163834

# Row 7
$ws.Cells.Item(7,1).Value2 = "MuSCs"
$ws.Cells.Item(7,2).Value2 = "Hp"
$ws.Cells.Item(7,3).Value2 = "Itgam"
$ws.Cells.Item(7,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(7,5).Value2 = 1
$ws.Cells.Item(7,6).Value2 = 0.3333333333333333
$ws.Cells.Item(7,7).Value2 = 0.1062546666666667
$ws.Cells.Item(7,8).Value2 = 0.318764
$ws.Cells.Item(7,9).Value2 = 0.01649928200866351
$ws.Cells.Item(7,10).Value2 = 0.01649928200866351
$ws.Cells.Item(7,11).Value2 = 3
$ws.Cells.Item(7,12).Value2 = 1
$ws.Cells.Item(7,13).Value2 = 38.416166
$ws.Cells.Item(7,14).Value2 = 115.248498
$ws.Cells.Item(7,15).Value2 = 0.9992910793627116
$ws.Cells.Item(7,16).Value2 = 0.9992910793627116
$ws.Cells.Item(7,17).Value2 = 4.081896912941333
$ws.Cells.Item(7,18).Value2 = 36.737072216472
$ws.Cells.Item(7,19).Value2 = 0.01648758532714713
$ws.Cells.Item(7,20).Value2 = 0.01648758532714713

# Row 8
$ws.Cells.Item(8,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(8,2).Value2 = "Hp"
$ws.Cells.Item(8,3).Value2 = "Itgam"
$ws.Cells.Item(8,4).Value2 = "MuSCs"
$ws.Cells.Item(8,5).Value2 = 2
$ws.Cells.Item(8,6).Value2 = 0.6666666666666666
$ws.Cells.Item(8,7).Value2 = 0.3992853333333333
$ws.Cells.Item(8,8).Value2 = 1.197856
$ws.Cells.Item(8,9).Value2 = 0.06200124214079897
$ws.Cells.Item(8,10).Value2 = 0.06200124214079897
$ws.Cells.Item(8,11).Value2 = 1
$ws.Cells.Item(8,12).Value2 = 0.3333333333333333
$ws.Cells.Item(8,13).Value2 = 0.02725333333333333
$ws.Cells.Item(8,14).Value2 = 0.08176
$ws.Cells.Item(8,15).Value2 = 0.0007089206372884383
$ws.Cells.Item(8,16).Value2 = 0.0007089206372884382
$ws.Cells.Item(8,17).Value2 = 0.01088185628444444
$ws.Cells.Item(8,18).Value2 = 0.09793670656
$ws.Cells.Item(8,19).Value2 = 0.00004395396009112999
$ws.Cells.Item(8,20).Value2 = 0.00004395396009112997

# Row 9
$ws.Cells.Item(9,1).Value2 = "Resolving-Mac"
$ws.Cells.Item(9,2).Value2 = "Hp"
$ws.Cells.Item(9,3).Value2 = "Itgam"
$ws.Cells.Item(9,4).Value2 = "Resolving-Mac"
$ws.Cells.Item(9,5).Value2 = 2
$ws.Cells.Item(9,6).Value2 = 0.6666666666666666
$ws.Cells.Item(9,7).Value2 = 0.3992853333333333
$ws.Cells.Item(9,8).Value2 = 1.197856
$ws.Cells.Item(9,9).Value2 = 0.06200124214079897
$ws.Cells.Item(9,10).Value2 = 0.06200124214079897
$ws.Cells.Item(9,11).Value2 = 3
$ws.Cells.Item(9,12).Value2 = 1
$ws.Cells.Item(9,13).Value2 = 38.416166
$ws.Cells.Item(9,14).Value2 = 115.248498
$ws.Cells.Item(9,15).Value2 = 0.9992910793627116
$ws.Cells.Item(9,16).Value2 = 0.9992910793627116
$ws.Cells.Item(9,17).Value2 = 15.33901164669867
$ws.Cells.Item(9,18).Value2 = 138.051104820288
$ws.Cells.Item(9,19).Value2 = 0.06195728818070784
$ws.Cells.Item(9,20).Value2 = 0.06195728818070784
